$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header cell in H1, copying the formatting of the existing
# header cells (bold, border, centered) from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new "Save" data values in column H for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
